# Apply the "[1.17.7]" version-history entry to the version log worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Push formatting down: the row that used to be the trailing blank
#        row (row 21) becomes the new data row, and a fresh blank row
#        (row 22), formatted like the old blank row, is appended below it.
$ws.Range("A21:C21").Copy() | Out-Null
$ws.Range("A22:C22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 21's Details/Date columns should look like the rest of the table
# (wrapped text column, left/top aligned date column).
$ws.Range("B20").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C20").Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Fill in the new version entry values.
$ws.Range("A21").Value = "[1.17.7]"

$detailsLines = @(
    "[handling bad templates - input file variants]",
    "- handle the variation of input files where we only have 5 address fields",
    "- adjust the error reported in the interface for the case of contacts with bad templates to show the First and Last names of the contacts as well as the opportunity",
    "- enable the tool to generate outputs when the custom scenario is being used, even if the contact has an invalid template assigned"
)
$details = [string]::Join("`n", $detailsLines)
$ws.Range("B21").Value = $details

$ws.Range("C21").Value = (Get-Date -Year 2018 -Month 9 -Day 16 -Hour 0 -Minute 0 -Second 0)

# Row 21 needs to grow to fit the wrapped text, same as other long entries.
$ws.Rows.Item(21).RowHeight = 90

# --- 3. Grow the table ("Table2") so it covers the new row too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C22"))

# --- 4. Update the view so the newly added row is visible/selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 20
$ws.Range("B21").Select() | Out-Null
